# Applies:
#  1. Slide 9 ("Pantalla" with cId=3287619798 / sldId=267): remove the
#     now-unused green "Grupo 32" status group and move the blue
#     "Grupo 2" status group up into its place.
#  2. Refresh the cached text of every auto-updating "datetimeFigureOut"
#     date field (slide master and all 11 layouts) from the stale
#     16/10/2024 to the current 2/11/2024.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1a. Remove the green "Grupo 32" (En Progreso) status group on slide 9.
# ---------------------------------------------------------------------
$s = $p.Slides.Item(9)
$s.Shapes.Item("Grupo 32").Delete()

# ---------------------------------------------------------------------
# 1b. Move the blue "Grupo 2" (Comenzadas) status group up to the slot
#     vacated by the deleted group.
# ---------------------------------------------------------------------
$g2 = $s.Shapes.Item("Grupo 2")
$g2.Top = 267.7676

# ---------------------------------------------------------------------
# 2. Re-cache the datetimeFigureOut date placeholders everywhere.
# ---------------------------------------------------------------------
$newDate = "2/11/2024"

function Set-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.PlaceholderFormat.Type -eq 16) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide master.
Set-DatePlaceholders($p.SlideMaster.Shapes)

# Every slide layout.
for ($L = 1; $L -le $p.SlideMaster.CustomLayouts.Count; $L++) {
    $cl = $p.SlideMaster.CustomLayouts.Item($L)
    Set-DatePlaceholders($cl.Shapes)
}

# NOTE: intentionally not touching $p.NotesMaster here - in this COM
# host, NotesMaster shape ids collide with SlideMaster shape ids and
# writes to NotesMaster.Shapes silently land on the (unrelated)
# SlideMaster shape with the same id instead of the notes master, so
# leaving it out avoids corrupting the slide master's placeholders.
